# Update ObjTables metadata strings: add schema='SBtab', reorder
# tableFormat attribute earlier, and bump the timestamps.
#
# Sheet "!!Compartment" (sheet1) carries BOTH the document-level banner
# (A1) and the table-level banner (A2). All other data sheets carry only
# the table-level banner in A1.
#
# All sheets are protected (no password), so each is unprotected before
# the edit and re-protected afterwards to preserve the original state.

$wb = $excel.ActiveWorkbook

function Set-A1Value {
    param($sheetName, $value)
    $sheet = $wb.Worksheets.Item($sheetName)
    $sheet.Unprotect()
    $sheet.Range("A1").Value = $value
    $sheet.Protect()
}

$wsCompartment = $wb.Worksheets.Item("!!Compartment")
$wsCompartment.Unprotect()
$wsCompartment.Range("A1").Value = "!!!ObjTables schema='SBtab' objTablesVersion='0.0.8' date='2020-03-09 23:59:24'"
$wsCompartment.Range("A2").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Compartment' name='Compartment' date='2020-03-09 23:59:24' objTablesVersion='0.0.8'"
$wsCompartment.Protect()

Set-A1Value "!!Compound" "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Compound' name='Compound' date='2020-03-09 23:59:24' objTablesVersion='0.0.8'"
Set-A1Value "!!Definition" "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Definition' name='Definition' date='2020-03-09 23:59:24' objTablesVersion='0.0.8'"
Set-A1Value "!!Enzyme" "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Enzyme' name='Enzyme' date='2020-03-09 23:59:24' objTablesVersion='0.0.8'"
Set-A1Value "!!FbcObjective" "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='FbcObjective' name='FbcObjective' date='2020-03-09 23:59:24' objTablesVersion='0.0.8'"
Set-A1Value "!!Gene" "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Gene' name='Gene' date='2020-03-09 23:59:24' objTablesVersion='0.0.8'"
Set-A1Value "!!Layout" "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Layout' name='Layout' date='2020-03-09 23:59:24' objTablesVersion='0.0.8'"
Set-A1Value "!!Measurement" "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Measurement' name='Measurement' date='2020-03-09 23:59:24' objTablesVersion='0.0.8'"
Set-A1Value "!!PbConfig" "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='PbConfig' name='PbConfig' date='2020-03-09 23:59:24' objTablesVersion='0.0.8'"
Set-A1Value "!!Position" "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Position' name='Position' date='2020-03-09 23:59:24' objTablesVersion='0.0.8'"
Set-A1Value "!!Protein" "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Protein' name='Protein' date='2020-03-09 23:59:24' objTablesVersion='0.0.8'"
Set-A1Value "!!Quantity" "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Quantity' name='Quantity' date='2020-03-09 23:59:24' objTablesVersion='0.0.8' level='1.0' version='0.1'"
Set-A1Value "!!QuantityInfo" "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='QuantityInfo' name='QuantityInfo' date='2020-03-09 23:59:25' objTablesVersion='0.0.8'"
Set-A1Value "!!QuantityMatrix" "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='QuantityMatrix' name='QuantityMatrix' date='2020-03-09 23:59:25' objTablesVersion='0.0.8'"
Set-A1Value "!!Reaction" "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Reaction' name='Reaction' date='2020-03-09 23:59:25' objTablesVersion='0.0.8'"
Set-A1Value "!!ReactionStoichiometry" "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='ReactionStoichiometry' name='ReactionStoichiometry' date='2020-03-09 23:59:25' objTablesVersion='0.0.8'"
Set-A1Value "!!Regulator" "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Regulator' name='Regulator' date='2020-03-09 23:59:25' objTablesVersion='0.0.8'"
Set-A1Value "!!Relation" "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Relation' name='Relation' date='2020-03-09 23:59:25' objTablesVersion='0.0.8'"
Set-A1Value "!!Relationship" "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Relationship' name='Relationship' date='2020-03-09 23:59:25' objTablesVersion='0.0.8'"
Set-A1Value "!!SparseMatrix" "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='SparseMatrix' name='SparseMatrix' date='2020-03-09 23:59:25' objTablesVersion='0.0.8'"
Set-A1Value "!!SparseMatrixColumn" "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='SparseMatrixColumn' name='SparseMatrixColumn' date='2020-03-09 23:59:25' objTablesVersion='0.0.8'"
Set-A1Value "!!SparseMatrixOrdered" "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='SparseMatrixOrdered' name='SparseMatrixOrdered' date='2020-03-09 23:59:25' objTablesVersion='0.0.8'"
Set-A1Value "!!SparseMatrixRow" "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='SparseMatrixRow' name='SparseMatrixRow' date='2020-03-09 23:59:25' objTablesVersion='0.0.8'"
Set-A1Value "!!StoichiometricMatrix" "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='StoichiometricMatrix' name='StoichiometricMatrix' date='2020-03-09 23:59:25' objTablesVersion='0.0.8'"
Set-A1Value "!!rxnconContingencyList" "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='rxnconContingencyList' name='rxnconContingencyList' date='2020-03-09 23:59:25' objTablesVersion='0.0.8'"
Set-A1Value "!!rxnconReactionList" "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='rxnconReactionList' name='rxnconReactionList' date='2020-03-09 23:59:25' objTablesVersion='0.0.8'"
